$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.383.86'
$ws.Range("E2").Value = '  +2.52%  '
$ws.Range("D3").Value = '3.591.34'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '243.23'
$ws.Range("E5").Value = '  +2.92%  '
$ws.Range("E6").Value = '  +17.25%  '
$ws.Range("D7").Value = '652.58'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").Value = '0.433'
$ws.Range("E8").Value = '  +8.50%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("D11").Value = '3.586.96'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '44.27'
$ws.Range("E12").Value = '  +4.39%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").Value = '4.260.61'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").Value = '97.159.47'
$ws.Range("E16").Value = '  +2.34%  '
$ws.Range("D17").Value = '0.0000262'
$ws.Range("E17").Value = '  +3.42%  '
$ws.Range("D18").Value = '8.67'
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("D19").Value = '3.573.46'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '12.56'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").Value = '18.09'
$ws.Range("E21").Value = '  +1.57%  '
$ws.Range("D22").Value = '0.528'
$ws.Range("E22").Value = '  +9.26%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '518.51'
$ws.Range("E23").Value = '  +1.96%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = '3.48'
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("D25").Value = '0.0000209'
$ws.Range("E25").Value = '  +6.53%  '
$ws.Range("D26").Value = '6.93'
$ws.Range("E26").Value = '  +1.42%  '
$ws.Range("D27").Value = '102.86'
$ws.Range("E27").Value = '  +8.17%  '
$ws.Range("E28").Value = '  +4.11%  '
$ws.Range("D29").Value = '3.784.96'
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").Value = '0.169'
$ws.Range("E30").Value = '  +17.13%  '
$ws.Range("D31").Value = '12.02'
$ws.Range("E31").Value = '  +4.15%  '
$ws.Range("D32").Value = '2.98'
$ws.Range("E32").Value = '  -1.76%  '
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("D34").Value = '0.187'
$ws.Range("E34").Value = '  +6.15%  '
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").Value = '31.88'
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = '0.575'
$ws.Range("E37").Value = '  +2.78%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '616.83'
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("D39").Value = '8.75'
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("E40").Value = '  -4.21%  '
$ws.Range("E41").Value = '  +2.68%  '
$ws.Range("E42").Value = '  +6.32%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").Value = '0.928'
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("D45").Value = '6.05'
$ws.Range("E45").Value = '  +4.78%  '
$ws.Range("D46").Value = '0.436'
$ws.Range("E46").Value = '  +40.65%  '
$ws.Range("D47").Value = '0.0444'
$ws.Range("E47").Value = '  +7.32%  '
$ws.Range("D48").Value = '2.32'
$ws.Range("E48").Value = '  +1.67%  '
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("D50").Value = '8.61'
$ws.Range("E50").Value = '  +5.24%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '3.30'
$ws.Range("E51").Value = '  +7.85%  '
